$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04217225448167028
$ws.Range("D2").Value = 0.6254007317907622
$ws.Range("E2").Value = 0.07971969209716079
$ws.Range("F2").Value = 7.67766579295477
$ws.Range("G2").Value = 0.002619586479155167
$ws.Range("I2").Value = 6.702089430277567
$ws.Range("L2").Value = 0.2309916117001336
$ws.Range("C3").Value = 0.03922197126684068
$ws.Range("D3").Value = 0.5997560311011512
$ws.Range("E3").Value = 0.07954608488602943
$ws.Range("F3").Value = 7.312664915336995
$ws.Range("G3").Value = 0.002637845245989812
$ws.Range("I3").Value = 6.407586212502139
$ws.Range("L3").Value = 0.2251548882320122
$ws.Range("C4").Value = 0.03740463833784702
$ws.Range("D4").Value = 0.584639479713104
$ws.Range("E4").Value = 0.07949401871503703
$ws.Range("F4").Value = 7.096297246036329
$ws.Range("G4").Value = 0.002649561795740255
$ws.Range("I4").Value = 6.233440295295793
$ws.Range("L4").Value = 0.2217824092075489
$ws.Range("C5").Value = 0.03666219987602659
$ws.Range("D5").Value = 0.5786322765800378
$ws.Range("E5").Value = 0.07948633921431458
$ws.Range("F5").Value = 7.010003513989318
$ws.Range("G5").Value = 0.002654464635297725
$ws.Range("I5").Value = 6.164096758339838
$ws.Range("L5").Value = 0.2204598454046476
$ws.Range("C6").Value = 0.03653879501806756
$ws.Range("D6").Value = 0.5776438682993899
$ws.Range("E6").Value = 0.07948587676149543
$ws.Range("F6").Value = 6.995785942668817
$ws.Range("G6").Value = 0.002655286526139592
$ws.Range("I6").Value = 6.15267866553549
$ws.Range("L6").Value = 0.2202433208500878
$ws.Range("C7").Value = 0.03739463359233497
$ws.Range("D7").Value = 0.5845578522429378
$ws.Range("E7").Value = 0.07949386055603647
$ws.Range("F7").Value = 7.095125945105792
$ws.Range("G7").Value = 0.002649627396396769
$ws.Range("I7").Value = 6.232498609579125
$ws.Range("L7").Value = 0.22176436490588
$ws.Range("C8").Value = 0.0411559705254092
$ws.Range("D8").Value = 0.6164245101484482
$ws.Range("E8").Value = 0.07964840728961775
$ws.Range("F8").Value = 7.55016332966801
$ws.Range("G8").Value = 0.002625777846402963
$ws.Range("I8").Value = 6.599122840620254
$ws.Range("L8").Value = 0.2289344036841214
$ws.Range("C9").Value = 0.04850408655065053
$ws.Range("D9").Value = 0.6841645191993848
$ws.Range("E9").Value = 0.08039250942737297
$ws.Range("F9").Value = 8.507258773144088
$ws.Range("G9").Value = 0.002582968601257977
$ws.Range("I9").Value = 7.373821485815768
$ws.Range("L9").Value = 0.244737444301208
$ws.Range("C10").Value = 0.05391293723054957
$ws.Range("D10").Value = 0.7375099410671737
$ws.Range("E10").Value = 0.08122049725153246
$ws.Range("F10").Value = 9.254813788462059
$ws.Range("G10").Value = 0.002553855825611224
$ws.Range("I10").Value = 7.981023838710541
$ws.Range("L10").Value = 0.2575078008063656
$ws.Range("C11").Value = 0.05638207931409056
$ws.Range("D11").Value = 0.7626428374326224
$ws.Range("E11").Value = 0.0816611850510327
$ws.Range("F11").Value = 9.605674321016181
$ws.Range("G11").Value = 0.002541102731789278
$ws.Range("I11").Value = 8.266458780448204
$ws.Range("L11").Value = 0.2635920434577343
$ws.Range("C12").Value = 0.05731886261362718
$ws.Range("D12").Value = 0.772292150533417
$ws.Range("E12").Value = 0.08183752697452107
$ws.Range("F12").Value = 9.740187184479396
$ws.Range("G12").Value = 0.002536342583986113
$ws.Range("I12").Value = 8.375952612406877
$ws.Range("L12").Value = 0.2659375074850487
$ws.Range("C13").Value = 0.05711702194680868
$ws.Range("D13").Value = 0.7702080016613877
$ws.Range("E13").Value = 0.08179912348012053
$ws.Range("F13").Value = 9.711142455711638
$ws.Range("G13").Value = 0.002537364712521625
$ws.Range("I13").Value = 8.35230729949069
$ws.Range("L13").Value = 0.2654304920336159
$ws.Range("C14").Value = 0.05645911036125995
$ws.Range("D14").Value = 0.7634340014603822
$ws.Range("E14").Value = 0.08167550163496173
$ws.Range("F14").Value = 9.616707137113508
$ws.Range("G14").Value = 0.002540709732774094
$ws.Range("I14").Value = 8.275438251062553
$ws.Range("L14").Value = 0.2637841629912145
$ws.Range("C15").Value = 0.05605636765984912
$ws.Range("D15").Value = 0.7593021559006843
$ws.Range("E15").Value = 0.08160101989789581
$ws.Range("F15").Value = 9.559080589523489
$ws.Range("G15").Value = 0.002542767622052709
$ws.Range("I15").Value = 8.228539289177945
$ws.Range("L15").Value = 0.2627812025925209
$ws.Range("C16").Value = 0.05375179145433151
$ws.Range("D16").Value = 0.7358855031885128
$ws.Range("E16").Value = 0.08119300721992317
$ws.Range("F16").Value = 9.232109634592632
$ws.Range("G16").Value = 0.002554699024685285
$ws.Range("I16").Value = 7.962562223481711
$ws.Range("L16").Value = 0.2571158785247718
$ws.Range("C17").Value = 0.05234058832979827
$ws.Range("D17").Value = 0.7217471205208881
$ws.Range("E17").Value = 0.08095928151066545
$ws.Range("F17").Value = 9.034355760268795
$ws.Range("G17").Value = 0.002562143222850413
$ws.Range("I17").Value = 7.801810180466532
$ws.Range("L17").Value = 0.2537121625092311
$ws.Range("C18").Value = 0.05152968032246008
$ws.Range("D18").Value = 0.7136961140038238
$ws.Range("E18").Value = 0.08083085946050517
$ws.Range("F18").Value = 8.921622795282815
$ws.Range("G18").Value = 0.002566471176759677
$ws.Range("I18").Value = 7.710212177574363
$ws.Range("L18").Value = 0.2517801903933616
$ws.Range("C19").Value = 0.05125523879461014
$ws.Range("D19").Value = 0.7109838738057022
$ws.Range("E19").Value = 0.08078840279664234
$ws.Range("F19").Value = 8.883623910112476
$ws.Range("G19").Value = 0.00256794453114891
$ws.Range("I19").Value = 7.679344388810819
$ws.Range("L19").Value = 0.2511304235253959
$ws.Range("C20").Value = 0.05249072895466611
$ws.Range("D20").Value = 0.7232437373293976
$ws.Range("E20").Value = 0.08098353811776704
$ws.Range("F20").Value = 9.055301795545347
$ws.Range("G20").Value = 0.002561345998619842
$ws.Range("I20").Value = 7.818832693362594
$ws.Range("L20").Value = 0.2540718144955036
$ws.Range("C21").Value = 0.05665230244720476
$ws.Range("D21").Value = 0.7654200440828163
$ws.Range("E21").Value = 0.0817115534302566
$ws.Range("F21").Value = 9.644399505822776
$ws.Range("G21").Value = 0.002539725353292699
$ws.Range("I21").Value = 8.297977735093468
$ws.Range("L21").Value = 0.2642665874597583
$ws.Range("C22").Value = 0.05938273211293676
$ws.Range("D22").Value = 0.7937580853726445
$ws.Range("E22").Value = 0.0822426413287296
$ws.Range("F22").Value = 10.03907326760577
$ws.Range("G22").Value = 0.002525997533138059
$ws.Range("I22").Value = 8.619361651941745
$ws.Range("L22").Value = 0.2711723072249299
$ws.Range("C23").Value = 0.05792429713754643
$ws.Range("D23").Value = 0.778560172045843
$ws.Range("E23").Value = 0.08195404232722936
$ws.Range("F23").Value = 9.827510510781678
$ws.Range("G23").Value = 0.002533287959481602
$ws.Range("I23").Value = 8.447051655679559
$ws.Range("L23").Value = 0.26746369757916
$ws.Range("C24").Value = 0.05242284913688877
$ws.Range("D24").Value = 0.7225668765433397
$ws.Range("E24").Value = 0.08097255319206553
$ws.Range("F24").Value = 9.045829115186393
$ws.Range("G24").Value = 0.00256170627373099
$ws.Range("I24").Value = 7.811134265689077
$ws.Range("L24").Value = 0.2539091385158372
$ws.Range("C25").Value = 0.04651647950163351
$ws.Range("D25").Value = 0.6652393358796473
$ws.Range("E25").Value = 0.08014275415760963
$ws.Range("F25").Value = 8.240914992244484
$ws.Range("G25").Value = 0.002594133311976418
$ws.Range("I25").Value = 7.157871378859795
$ws.Range("L25").Value = 0.2402659826008886
